# Update automàtic: dades i banners [2026-02-20 09:56]
# Update DATA_EXTRACCIO timestamps in column H of "Dades_Període" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Període")

$ws.Range("H2").Value = "2026-02-20 09:56:47"
$ws.Range("H3").Value = "2026-02-20 09:56:49"
$ws.Range("H4").Value = "2026-02-20 09:56:49"
$ws.Range("H5").Value = "2026-02-20 09:56:49"
$ws.Range("H6").Value = "2026-02-20 09:56:49"
